$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.91
$ws.Range("I2").Value = 4.33
$ws.Range("Y2").Value = 6.5
$ws.Range("Z2").Value = 8.5
$ws.Range("AA2").Value = 8.5
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 15
$ws.Range("AD2").Value = 29
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 21
$ws.Range("AL2").Value = 15
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 34
$ws.Range("AO2").Value = 41

# Row 3
$ws.Range("Q3").Value = 1.29
$ws.Range("R3").Value = 3.75
$ws.Range("Z3").Value = 9
$ws.Range("AI3").Value = 151
$ws.Range("AL3").Value = 34

# Row 12
$ws.Range("G12").Value = 2.15
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 3.3
$ws.Range("J12").Value = 2.67
$ws.Range("K12").Value = 2.07
$ws.Range("L12").Value = 3.8
$ws.Range("O12").Value = 1.34
$ws.Range("P12").Value = 2.75
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.65
$ws.Range("S12").Value = 3.25
$ws.Range("T12").Value = 1.25
$ws.Range("U12").Value = 1.39
$ws.Range("V12").Value = 2.57
$ws.Range("Y12").Value = 7.1
$ws.Range("Z12").Value = 10.25
$ws.Range("AA12").Value = 8.75
$ws.Range("AB12").Value = 21
$ws.Range("AC12").Value = 18
$ws.Range("AD12").Value = 29
$ws.Range("AE12").Value = 8.5
$ws.Range("AF12").Value = 6.1
$ws.Range("AJ12").Value = 9
$ws.Range("AK12").Value = 17
$ws.Range("AL12").Value = 11.75
$ws.Range("AM12").Value = 45
$ws.Range("AN12").Value = 32
$ws.Range("AO12").Value = 40

# Row 14
$ws.Range("Q14").Value = 1.9
$ws.Range("R14").Value = 1.95
$ws.Range("S14").Value = 3.25
$ws.Range("T14").Value = 1.33

# Row 18
$ws.Range("G18").Value = 2.4
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 3.2
$ws.Range("J18").Value = 3.25
$ws.Range("K18").Value = 1.83
$ws.Range("L18").Value = 4
$ws.Range("Y18").Value = 6
$ws.Range("Z18").Value = 10
$ws.Range("AA18").Value = 11
$ws.Range("AB18").Value = 23
$ws.Range("AJ18").Value = 7
$ws.Range("AK18").Value = 13
$ws.Range("AM18").Value = 34
$ws.Range("AO18").Value = 41

# Row 30
$ws.Range("G30").Value = 3.2
$ws.Range("I30").Value = 2.1
$ws.Range("J30").Value = 3.75
$ws.Range("AA30").Value = 12
$ws.Range("AJ30").Value = 9
$ws.Range("AM30").Value = 19
$ws.Range("AN30").Value = 15

# Row 31
$ws.Range("AQ31").Value = 1.36
$ws.Range("AS31").Value = 2.29

# Row 32
$ws.Range("Q32").Value = 1.9
$ws.Range("R32").Value = 1.95

# Row 34
$ws.Range("G34").Value = 2.3
$ws.Range("I34").Value = 3.4
$ws.Range("J34").Value = 3.2
$ws.Range("L34").Value = 4.33
$ws.Range("M34").Value = 1.1
$ws.Range("N34").Value = 7
$ws.Range("Y34").Value = 6
$ws.Range("Z34").Value = 9.5
$ws.Range("AB34").Value = 21
$ws.Range("AJ34").Value = 8
$ws.Range("AK34").Value = 15
$ws.Range("AL34").Value = 13
$ws.Range("AM34").Value = 41
$ws.Range("AN34").Value = 34

# Row 37
$ws.Range("G37").Value = 4.55
$ws.Range("H37").Value = 3
$ws.Range("I37").Value = 1.82
$ws.Range("J37").Value = 5.2
$ws.Range("L37").Value = 2.45
$ws.Range("N37").Value = 5.5
$ws.Range("S37").Value = 4.3
$ws.Range("W37").Value = 2.15
$ws.Range("Z37").Value = 24
$ws.Range("AA37").Value = 16.5
$ws.Range("AD37").Value = 80
$ws.Range("AE37").Value = 5.5
$ws.Range("AH37").Value = 120
$ws.Range("AJ37").Value = 5.3
$ws.Range("AK37").Value = 7.5
$ws.Range("AN37").Value = 17.5

# Row 38
$ws.Range("G38").Value = 2.27
$ws.Range("H38").Value = 3.05
$ws.Range("I38").Value = 3.2
$ws.Range("J38").Value = 2.82
$ws.Range("N38").Value = 6.5
$ws.Range("P38").Value = 2.92
$ws.Range("S38").Value = 3.45
$ws.Range("W38").Value = 1.78
$ws.Range("X38").Value = 1.93
$ws.Range("Y38").Value = 7.3
$ws.Range("Z38").Value = 10.75
$ws.Range("AA38").Value = 8.75
$ws.Range("AB38").Value = 23
$ws.Range("AE38").Value = 6.5
$ws.Range("AF38").Value = 5.9
$ws.Range("AJ38").Value = 9
$ws.Range("AL38").Value = 11
$ws.Range("AN38").Value = 30

# Row 45
$ws.Range("M45").Value = 1.05
$ws.Range("N45").Value = 11

# Row 46
$ws.Range("H46").Value = 3
$ws.Range("I46").Value = 2.67
$ws.Range("W46").Value = 1.78
$ws.Range("X46").Value = 1.82
$ws.Range("Y46").Value = 7.6
$ws.Range("AA46").Value = 9.75
$ws.Range("AE46").Value = 7.8
$ws.Range("AF46").Value = 5.8
$ws.Range("AK46").Value = 13
$ws.Range("AM46").Value = 32
